# Auto-generated Excel COM-interop script to apply the diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 199
$ws.Range("I12").Value = 199
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 199
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -29
$ws.Range("N12").ClearContents()
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 2301
$ws.Range("I29").Value = 2301
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 6903
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -6622
$ws.Range("N29").ClearContents()
$ws.Range("H32").Value = 6116.4
$ws.Range("I32").Value = 6966.6665
$ws.Range("J32").Value = 4841
$ws.Range("K32").Value = 6966.6665
$ws.Range("L32").Value = 4841
$ws.Range("M32").Value = -6640.6665
$ws.Range("N32").Value = -5493
$ws.Range("H58").Value = 8502.333000000001
$ws.Range("I58").Value = 2500
$ws.Range("J58").Value = 9702.799999999999
$ws.Range("K58").Value = 7500
$ws.Range("L58").Value = 29108.4
$ws.Range("M58").Value = -7350
$ws.Range("N58").Value = -29408.4
$ws.Range("H137").Value = 1945.6316
$ws.Range("I137").Value = 1305.3846
$ws.Range("K137").Value = 3916.1538
$ws.Range("M137").Value = -1366.1538
$ws.Range("H138").Value = 7512.4814
$ws.Range("J138").Value = 7972.625
$ws.Range("L138").Value = 23917.875
$ws.Range("N138").Value = -34197.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H16").Value = 655
$ws.Range("I16").Value = 655
$ws.Range("K16").Value = 655
$ws.Range("M16").Value = -368
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H45").Value = 3266.3333
$ws.Range("I45").Value = 2822.6
$ws.Range("J45").Value = 3821
$ws.Range("K45").Value = 2822.6
$ws.Range("L45").Value = 3821
$ws.Range("M45").Value = -2445.6
$ws.Range("N45").Value = -4575
$ws.Range("H74").Value = 1763.8572
$ws.Range("I74").Value = 1057.8334
$ws.Range("K74").Value = 1057.8334
$ws.Range("M74").Value = -183.8334
$ws.Range("H77").Value = 1763.8572
$ws.Range("I77").Value = 1057.8334
$ws.Range("K77").Value = 5289.166999999999
$ws.Range("M77").Value = -921.1669999999995
$ws.Range("H110").Value = 3628.5386
$ws.Range("I110").Value = 3727.4
$ws.Range("K110").Value = 3727.4
$ws.Range("M110").Value = -1682.4
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 3955.1428
$ws.Range("I132").Value = 3475.6
$ws.Range("K132").Value = 10426.8
$ws.Range("M132").Value = -7896.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H19").Value = 4000
$ws.Range("J19").Value = 4000
$ws.Range("L19").Value = 4000
$ws.Range("N19").Value = -4346
$ws.Range("H107").Value = 2811
$ws.Range("I107").Value = 2811
$ws.Range("K107").Value = 2811
$ws.Range("M107").Value = -891
$ws.Range("H134").Value = 4211.5
$ws.Range("I134").Value = 4211.5
$ws.Range("K134").Value = 12634.5
$ws.Range("M134").Value = -10099.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H31").Value = 3933.3157
$ws.Range("I31").Value = 1541.7693
$ws.Range("K31").Value = 1541.7693
$ws.Range("M31").Value = -1246.7693
$ws.Range("H34").Value = 3933.3157
$ws.Range("I34").Value = 1541.7693
$ws.Range("K34").Value = 1541.7693
$ws.Range("M34").Value = -1339.7693
$ws.Range("H58").Value = 2415.2
$ws.Range("I58").Value = 2420.25
$ws.Range("K58").Value = 2420.25
$ws.Range("M58").Value = -2217.25
$ws.Range("H62").Value = 2367.5
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 2205.7144
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 2205.7144
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -3453.7144
$ws.Range("H65").Value = 2367.5
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 2205.7144
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 11028.572
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -17268.572
$ws.Range("H134").Value = 13448
$ws.Range("I134").Value = 13448
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 40344
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -37809
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2415.2
$ws.Range("I136").Value = 2420.25
$ws.Range("K136").Value = 7260.75
$ws.Range("M136").Value = -4710.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 194489
$ws.Range("J37").Value = 194489
$ws.Range("L37").Value = 583467
$ws.Range("N37").Value = -583691

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 6333.3335
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 6333.3335
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 6333.3335
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -6565.3335
$ws.Range("H132").Value = 5117
$ws.Range("I132").Value = 4842.5557
$ws.Range("J132").Value = 5734.5
$ws.Range("K132").Value = 14527.6671
$ws.Range("L132").Value = 17203.5
$ws.Range("M132").Value = -11997.6671
$ws.Range("N132").Value = -22263.5
$ws.Range("H135").Value = 125000
$ws.Range("J135").Value = 125000
$ws.Range("L135").Value = 125000
$ws.Range("N135").Value = -135140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 471
$ws.Range("I55").Value = 414.66666
$ws.Range("J55").Value = 640
$ws.Range("K55").Value = 414.66666
$ws.Range("L55").Value = 640
$ws.Range("M55").Value = -241.66666
$ws.Range("N55").Value = -986
$ws.Range("H61").Value = 2600
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4798
$ws.Range("H113").Value = 2600
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2286.4055
$ws.Range("I132").Value = 1962.5862
$ws.Range("J132").Value = 3460.25
$ws.Range("K132").Value = 5887.7586
$ws.Range("L132").Value = 10380.75
$ws.Range("M132").Value = -3357.7586
$ws.Range("N132").Value = -15440.75
